$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Wnt5a", "Fzd2", "ECs", 2, 1, 0.125615, 0.25123, 0.02647478672532295, 0.01780700335556722, 2, 1, 0.1557005, 0.311401, 0.02102398211576467, 0.01500040222529337, 0.0195583183075, 0.07823327323, 0.0005566054426318736, 0.000267112212760657),
    @("ECs", "Wnt5a", "Fzd2", "FAPs", 2, 1, 0.125615, 0.25123, 0.02647478672532295, 0.01780700335556722, 3, 1, 5.923689, 17.771067, 0.7998659708565604, 0.8560446272575798, 0.744104193735, 4.464625162410001, 0.02117628098727082, 0.01524358955009101),
    @("ECs", "Wnt5a", "Fzd2", "Inflammatory-Mac", 2, 1, 0.125615, 0.25123, 0.02647478672532295, 0.01780700335556722, 1, 0.3333333333333333, 0.024117, 0.072351, 0.003256478795417461, 0.003485197868350457, 0.003029456955, 0.01817674173, 0.00008621458158421389, 0.0000620609301365323),
    @("ECs", "Wnt5a", "Fzd2", "MuSCs", 2, 1, 0.125615, 0.25123, 0.02647478672532295, 0.01780700335556722, 2, 1, 1.3023455, 2.604691, 0.1758535682322574, 0.1254697726487764, 0.1635941299825, 0.65437651993, 0.004655685713836042, 0.002234240662579017),
    @("FAPs", "Wnt5a", "Fzd2", "ECs", 3, 1, 4.619088000000001, 13.857264, 0.9735252132746771, 0.9821929966444328, 2, 1, 0.1557005, 0.311401, 0.02102398211576467, 0.01500040222529337, 0.719194311144, 4.315165866864, 0.02046737667313279, 0.01473329001253271),
    @("FAPs", "Wnt5a", "Fzd2", "FAPs", 3, 1, 4.619088000000001, 13.857264, 0.9735252132746771, 0.9821929966444328, 3, 1, 5.923689, 17.771067, 0.7998659708565604, 0.8560446272575798, 27.362040775632, 246.258366980688, 0.7786896898692897, 0.8408010377074888),
    @("FAPs", "Wnt5a", "Fzd2", "Inflammatory-Mac", 3, 1, 4.619088000000001, 13.857264, 0.9735252132746771, 0.9821929966444328, 1, 0.3333333333333333, 0.024117, 0.072351, 0.003256478795417461, 0.003485197868350457, 0.111398545296, 1.002586907664, 0.003170264213833247, 0.003423136938213925),
    @("FAPs", "Wnt5a", "Fzd2", "MuSCs", 3, 1, 4.619088000000001, 13.857264, 0.9735252132746771, 0.9821929966444328, 2, 1, 1.3023455, 2.604691, 0.1758535682322574, 0.1254697726487764, 6.015648470904001, 36.093890825424, 0.1711978825184213, 0.1232355319861974),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
